# ============================================================
# feat: add 2022-Q1 data
#
# The sheet that used to be named "总计" becomes the new
# "2022-Q1" sheet (reusing its sheetId/relationship id) and is
# repopulated with per-fund holdings for 2022-Q1. A brand-new
# "总计" sheet is appended after it, carrying the historical
# roll-up table plus the new 2022-Q1 summary row.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- Template cells to copy cell formatting (style index 2:
#      bold font + thin box border + center/top alignment) from.
$fmtSrc = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------
# 1) Turn the existing "总计" sheet into "2022-Q1" and replace
#    its contents with the per-fund holdings table.
# ------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Name = "2022-Q1"
$ws.Cells.Clear()

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"


# -- per-fund data rows (A: running index, H: rank -> plain numbers;
#    B/D/E/F/G carry a leading "'" so Excel stores them as text, just
#    like the source "基金代码/基金规模/..." columns on the other
#    per-quarter sheets; C is plain text; G10 is a genuine 0 number)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'004350"
$ws.Range("C2").Value = '汇丰晋信价值先锋股票'
$ws.Range("D2").Value = "'4.99"
$ws.Range("E2").Value = "'93.32"
$ws.Range("F2").Value = "'2.85"
$ws.Range("G2").Value = "'0.1422"
$ws.Range("H2").Value = 7
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'004206"
$ws.Range("C3").Value = '华商元亨灵活配置混合'
$ws.Range("D3").Value = "'5.64"
$ws.Range("E3").Value = "'29.94"
$ws.Range("F3").Value = "'1.63"
$ws.Range("G3").Value = "'0.0919"
$ws.Range("H3").Value = 6
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'290012"
$ws.Range("C4").Value = '泰信行业精选灵活配置混合A'
$ws.Range("D4").Value = "'0.76"
$ws.Range("E4").Value = "'92.62"
$ws.Range("F4").Value = "'6.55"
$ws.Range("G4").Value = "'0.0498"
$ws.Range("H4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'002681"
$ws.Range("C5").Value = '金鹰元和灵活配置混合A'
$ws.Range("D5").Value = "'0.56"
$ws.Range("E5").Value = "'81.63"
$ws.Range("F5").Value = "'4.15"
$ws.Range("G5").Value = "'0.0232"
$ws.Range("H5").Value = 7
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'002682"
$ws.Range("C6").Value = '金鹰元和灵活配置混合C'
$ws.Range("D6").Value = "'0.25"
$ws.Range("E6").Value = "'81.63"
$ws.Range("F6").Value = "'4.15"
$ws.Range("G6").Value = "'0.0104"
$ws.Range("H6").Value = 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'008135"
$ws.Range("C7").Value = '华宸未来价值先锋混合'
$ws.Range("D7").Value = "'0.20"
$ws.Range("E7").Value = "'86.99"
$ws.Range("F7").Value = "'4.07"
$ws.Range("G7").Value = "'0.0081"
$ws.Range("H7").Value = 7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'001448"
$ws.Range("C8").Value = '华商双翼平衡混合'
$ws.Range("D8").Value = "'0.38"
$ws.Range("E8").Value = "'39.74"
$ws.Range("F8").Value = "'1.97"
$ws.Range("G8").Value = "'0.0075"
$ws.Range("H8").Value = 7
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'004456"
$ws.Range("C9").Value = '兴银消费新趋势灵活配置混合'
$ws.Range("D9").Value = "'0.06"
$ws.Range("E9").Value = "'82.19"
$ws.Range("F9").Value = "'4.54"
$ws.Range("G9").Value = "'0.0027"
$ws.Range("H9").Value = 9
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'002583"
$ws.Range("C10").Value = '泰信行业精选灵活配置混合C'
$ws.Range("D10").Value = "'0.00"
$ws.Range("E10").Value = "'92.62"
$ws.Range("F10").Value = "'6.55"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 2

# -- apply the bold/boxed "header" style (same style index used by
#    every other sheet's header row + index column) to the new cells
$fmtSrc.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------
# 2) Append a brand-new "总计" sheet after "2022-Q1" holding the
#    historical roll-up (unchanged rows, shifted down by one)
#    plus a new row for 2022-Q1 at the top of the data.
# ------------------------------------------------------------
$newSheetTmp = $wb.Worksheets.Add()
$newSheetTmp.Name = "总计"
$newSheetTmp.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# NOTE: re-fetch by name after Move() -- the old object reference goes
# stale (rebinds to whatever sheet now sits at the pre-move index) once
# the tab order changes underneath it.
$total = $wb.Worksheets.Item("总计")

$total.PageSetup.LeftMargin = 0.75 * 72
$total.PageSetup.RightMargin = 0.75 * 72
$total.PageSetup.TopMargin = 1 * 72
$total.PageSetup.BottomMargin = 1 * 72
$total.PageSetup.HeaderMargin = 0.5 * 72
$total.PageSetup.FooterMargin = 0.5 * 72

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = '2022-Q1'
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.34
$total.Range("A3").Value = 1
$total.Range("B3").Value = '2021-Q4'
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.03
$total.Range("A4").Value = 2
$total.Range("B4").Value = '2021-Q3'
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.03
$total.Range("A5").Value = 3
$total.Range("B5").Value = '2021-Q2'
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.18
$total.Range("A6").Value = 4
$total.Range("B6").Value = '2021-Q1'
$total.Range("C6").Value = 11
$total.Range("D6").Value = 0.85
$total.Range("A7").Value = 5
$total.Range("B7").Value = '2020-Q4'
$total.Range("C7").Value = 5
$total.Range("D7").Value = 1.02

$fmtSrc.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
